$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing rows 1-27 shift down to 2-28,
# carrying their values/styles/formatting along with them.
$ws.Rows.Item(1).Insert()

# New row 1: sequential integers 0..11, keep header's bold/centered style
# by copying former header row's formatting (now row 2) onto row 1.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)

for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value2 = $c - 1
}

# Former header row (now row 2) should look like a normal data row:
# strip the header style/formatting, keep its text values.
$ws.Range("A2:L2").ClearFormats()

# I2 (blank in original header), K2 and L2 become empty for this row
$ws.Cells.Item(2, 9).Value2 = ""
$ws.Cells.Item(2, 11).Value2 = ""
$ws.Cells.Item(2, 12).Value2 = ""
